# Insert a new data row at row 617 (pushing the existing rows 617-646 down
# to 618-647) and populate it with the new "Macroferia Regional de Talca -
# Zanahoria" weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(617).Insert()

$ws.Range("A617").Value = 5
$ws.Range("B617").Value = "Macroferia Regional de Talca"
$ws.Range("C617").Value = "Maule"
$ws.Range("D617").Value = 45267
$ws.Range("E617").Value = 7
$ws.Range("F617").Value = 100114013
$ws.Range("G617").Value = "Zanahoria"
$ws.Range("H617").Value = "Sin especificar"
$ws.Range("I617").Value = "Primera"
$ws.Range("J617").Value = 500
$ws.Range("K617").Value = 5000
$ws.Range("L617").Value = 5000
$ws.Range("M617").Value = 5000
$ws.Range("N617").Value = "$/saco 20 kilos"
$ws.Range("O617").Value = "Provincia de Melipilla"
$ws.Range("P617").Value = 250
$ws.Range("Q617").Value = 20
$ws.Range("R617").Value = "Hortaliza"
